# Auto-generated Excel COM-interop script to apply the diff changes
# to Sheets/Sargatanas_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 3025.25
$ws.Range("I12").Value = 2149.5
$ws.Range("K12").Value = 2149.5
$ws.Range("M12").Value = -1979.5
$ws.Range("H17").Value = 1337.8485
$ws.Range("J17").Value = 1433.8928
$ws.Range("L17").Value = 4301.678400000001
$ws.Range("N17").Value = -4637.678400000001
$ws.Range("H18").Value = 4985.4
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H41").Value = 10417184
$ws.Range("I41").Value = 15625290
$ws.Range("J41").Value = 974.5
$ws.Range("K41").Value = 15625290
$ws.Range("L41").Value = 974.5
$ws.Range("M41").Value = -15624850
$ws.Range("N41").Value = -1854.5
$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 3666.6667
$ws.Range("L46").Value = 11000.0001
$ws.Range("N46").Value = -11238.0001
$ws.Range("H60").Value = 4000
$ws.Range("J60").Value = 3666.6667
$ws.Range("L60").Value = 11000.0001
$ws.Range("N60").Value = -11968.0001
$ws.Range("H116").Value = 22733180
$ws.Range("I116").Value = 35716824
$ws.Range("K116").Value = 35716824
$ws.Range("M116").Value = -35713382
$ws.Range("H135").Value = 323079.8
$ws.Range("I135").Value = 323079.8
$ws.Range("K135").Value = 2907718.2
$ws.Range("M135").Value = -2905183.2
$ws.Range("H137").Value = 2356.3333
$ws.Range("I137").Value = 2319.1667
$ws.Range("K137").Value = 6957.500100000001
$ws.Range("M137").Value = -4407.500100000001
$ws.Range("H138").Value = 4689.537
$ws.Range("I138").Value = 1732.9286
$ws.Range("J138").Value = 5470.5283
$ws.Range("K138").Value = 5198.7858
$ws.Range("L138").Value = 16411.5849
$ws.Range("M138").Value = -58.78579999999965
$ws.Range("N138").Value = -26691.5849
$ws.Range("H141").Value = 3710.077
$ws.Range("I141").Value = 3623.7
$ws.Range("K141").Value = 10871.1
$ws.Range("M141").Value = -5691.099999999999
$ws.Range("N18").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 57978.453
$ws.Range("I74").Value = 121269.8
$ws.Range("J74").Value = 5235.6665
$ws.Range("K74").Value = 121269.8
$ws.Range("L74").Value = 5235.6665
$ws.Range("M74").Value = -120395.8
$ws.Range("N74").Value = -6983.6665
$ws.Range("H77").Value = 57978.453
$ws.Range("I77").Value = 121269.8
$ws.Range("J77").Value = 5235.6665
$ws.Range("K77").Value = 606349
$ws.Range("L77").Value = 26178.3325
$ws.Range("M77").Value = -601981
$ws.Range("N77").Value = -34914.3325
$ws.Range("H102").Value = 789.55554
$ws.Range("I102").Value = 757.8261
$ws.Range("K102").Value = 757.8261
$ws.Range("M102").Value = 864.1739
$ws.Range("H110").Value = 41668044
$ws.Range("I110").Value = 1166.6666
$ws.Range("J110").Value = 166668670
$ws.Range("K110").Value = 1166.6666
$ws.Range("L110").Value = 166668670
$ws.Range("M110").Value = 878.3334
$ws.Range("N110").Value = -166672760
$ws.Range("H132").Value = 8376.468999999999
$ws.Range("I132").Value = 8735.75
$ws.Range("J132").Value = 8160.9
$ws.Range("K132").Value = 26207.25
$ws.Range("L132").Value = 24482.7
$ws.Range("M132").Value = -23677.25
$ws.Range("N132").Value = -29542.7

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3098.6
$ws.Range("I105").Value = 2456.4167
$ws.Range("J105").Value = 4061.875
$ws.Range("K105").Value = 2456.4167
$ws.Range("L105").Value = 4061.875
$ws.Range("M105").Value = -709.4167000000002
$ws.Range("N105").Value = -7555.875
$ws.Range("H107").Value = 48917620
$ws.Range("I107").Value = 70315790
$ws.Range("K107").Value = 70315790
$ws.Range("M107").Value = -70313870
$ws.Range("H134").Value = 6957.7896
$ws.Range("I134").Value = 2642.5
$ws.Range("K134").Value = 7927.5
$ws.Range("M134").Value = -5392.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5299.7617
$ws.Range("I16").Value = 3640.5
$ws.Range("J16").Value = 6808.1816
$ws.Range("K16").Value = 3640.5
$ws.Range("L16").Value = 6808.1816
$ws.Range("M16").Value = -3353.5
$ws.Range("N16").Value = -7382.1816
$ws.Range("H113").Value = 5299.7617
$ws.Range("I113").Value = 3640.5
$ws.Range("J113").Value = 6808.1816
$ws.Range("K113").Value = 3640.5
$ws.Range("L113").Value = 6808.1816
$ws.Range("M113").Value = -1470.5
$ws.Range("N113").Value = -11148.1816
$ws.Range("H132").Value = 5100.354
$ws.Range("I132").Value = 2391.1853
$ws.Range("J132").Value = 8583.571
$ws.Range("K132").Value = 7173.5559
$ws.Range("L132").Value = 25750.713
$ws.Range("M132").Value = -4643.5559
$ws.Range("N132").Value = -30810.713
$ws.Range("H134").Value = 4916.068
$ws.Range("I134").Value = 2143.9429
$ws.Range("K134").Value = 6431.8287
$ws.Range("M134").Value = -3896.8287

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3125336.8
$ws.Range("I12").Value = 975.75
$ws.Range("K12").Value = 2927.25
$ws.Range("M12").Value = -2754.25
$ws.Range("H131").Value = 1760.1111
$ws.Range("I131").Value = 1512.8334
$ws.Range("J131").Value = 2254.6667
$ws.Range("K131").Value = 4538.5002
$ws.Range("L131").Value = 6764.000100000001
$ws.Range("M131").Value = 501.4997999999996
$ws.Range("N131").Value = -16844.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.55556
$ws.Range("I2").Value = 63
$ws.Range("K2").Value = 63
$ws.Range("M2").Value = 50
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21136
$ws.Range("H97").Value = 1440.7407
$ws.Range("I97").Value = 1393.05
$ws.Range("K97").Value = 1393.05
$ws.Range("M97").Value = -897.05
$ws.Range("H132").Value = 4954.1562
$ws.Range("I132").Value = 1486.4783
$ws.Range("J132").Value = 13816
$ws.Range("K132").Value = 4459.4349
$ws.Range("L132").Value = 41448
$ws.Range("M132").Value = -1929.4349
$ws.Range("N132").Value = -46508
$ws.Range("H140").Value = 84389.5
$ws.Range("J140").Value = 84389.5
$ws.Range("L140").Value = 84389.5
$ws.Range("N140").Value = -94749.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1346.5
$ws.Range("I16").Value = 1346.5
$ws.Range("K16").Value = 1346.5
$ws.Range("M16").Value = -1176.5
$ws.Range("H21").Value = 7000
$ws.Range("J21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("N21").Value = -7348
$ws.Range("H122").Value = 4738.511
$ws.Range("I122").Value = 3517.2307
$ws.Range("J122").Value = 6409.737
$ws.Range("K122").Value = 10551.6921
$ws.Range("L122").Value = 19229.211
$ws.Range("M122").Value = -8101.6921
$ws.Range("N122").Value = -24129.211
$ws.Range("H132").Value = 8069854
$ws.Range("I132").Value = 14708164
$ws.Range("J132").Value = 9049.25
$ws.Range("K132").Value = 44124492
$ws.Range("L132").Value = 27147.75
$ws.Range("M132").Value = -44121962
$ws.Range("N132").Value = -32207.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 2824.75
$ws.Range("J41").Value = 2824.75
$ws.Range("L41").Value = 2824.75
$ws.Range("N41").Value = -3604.75
$ws.Range("H104").Value = 49235
$ws.Range("J104").Value = 49235
$ws.Range("L104").Value = 49235
$ws.Range("N104").Value = -56223
$ws.Range("H107").Value = 1313.6428
$ws.Range("I107").Value = 1398.5714
$ws.Range("K107").Value = 4195.7142
$ws.Range("M107").Value = -2275.7142
